$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 444, shifting rows 444:553 down to 445:554.
$ws.Rows.Item(444).Insert()

# Populate the new row 444 with a new data record (same static fields as the
# old row 444 carried down to row 445, but with updated Fecha / Volumen /
# Precio minimo / Precio maximo / Precio promedio ponderado / Origen /
# Precio $/Kg values).
$ws.Range("A444").Value = 7
$ws.Range("B444").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C444").Value = "Ñuble"
$ws.Range("D444").Value = (Get-Date -Year 2023 -Month 11 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E444").Value = 16
$ws.Range("F444").Value = 100114013
$ws.Range("G444").Value = "Zanahoria"
$ws.Range("H444").Value = "Sin especificar"
$ws.Range("I444").Value = "Primera"
$ws.Range("J444").Value = 200
$ws.Range("K444").Value = 6000
$ws.Range("L444").Value = 6000
$ws.Range("M444").Value = 6000
$ws.Range("N444").Value = "$/saco 20 kilos"
$ws.Range("O444").Value = "Provincia de Diguillín"
$ws.Range("P444").Value = 300
$ws.Range("Q444").Value = 20
$ws.Range("R444").Value = "Hortaliza"
